$d = $word.ActiveDocument

# Update the date line at the top of the document
$d.Content.Find.Execute("2024-06-25 Tuesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-06-26 Wednesday", 2)

# Update the division problems in the table. The table has data on rows
# 1, 5, 9, 13, 17 (1-based), each with 5 columns. Addressing each cell
# explicitly avoids any ambiguity from repeated values (e.g. "983÷5=").
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "646÷9="
$t.Cell(1,2).Range.Text = "744÷3="
$t.Cell(1,3).Range.Text = "877÷5="
$t.Cell(1,4).Range.Text = "936÷8="
$t.Cell(1,5).Range.Text = "730÷5="

$t.Cell(5,1).Range.Text = "683÷6="
$t.Cell(5,2).Range.Text = "654÷3="
$t.Cell(5,3).Range.Text = "931÷4="
$t.Cell(5,4).Range.Text = "409÷2="
$t.Cell(5,5).Range.Text = "348÷8="

$t.Cell(9,1).Range.Text = "470÷6="
$t.Cell(9,2).Range.Text = "471÷8="
$t.Cell(9,3).Range.Text = "872÷9="
$t.Cell(9,4).Range.Text = "651÷8="
$t.Cell(9,5).Range.Text = "951÷3="

$t.Cell(13,1).Range.Text = "616÷6="
$t.Cell(13,2).Range.Text = "274÷5="
$t.Cell(13,3).Range.Text = "709÷8="
$t.Cell(13,4).Range.Text = "983÷5="
$t.Cell(13,5).Range.Text = "722÷4="

$t.Cell(17,1).Range.Text = "845÷6="
$t.Cell(17,2).Range.Text = "232÷4="
$t.Cell(17,3).Range.Text = "665÷4="
$t.Cell(17,4).Range.Text = "165÷3="
$t.Cell(17,5).Range.Text = "195÷9="
